$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths -----------------------------------------------------
# Column D: 19.5703125 -> 23
$ws.Columns("D").ColumnWidth = 22.16
# New column O width -> 12.85546875 (closest achievable)
$ws.Columns("O").ColumnWidth = 12

# --- New number-format style (numFmtId 49 "@", Text) applied to D & O --
# All D value-cells except the formula cells in rows 86/87/92/93 and the
# text "separator" cells D85/D91, which are not touched by the edit
$ws.Range("D1:D84").NumberFormat = "@"
# All O cells (1-127)
$ws.Range("O1:O127").NumberFormat = "@"

# --- Formula changes -----------------------------------------------------
$ws.Range("D86").Formula = "=(82*2)+1"
$ws.Range("D92").Formula = "=(127*2)+1"

# --- Selection / view ----------------------------------------------------
$ws.Range("D5").Select()
